# Apply the "0.1.3" guide update:
#  - Metadata sheet: bump Version value from 0.1.2 to 0.1.3
#  - "Include from ISO 3166-1 Codes" sheet: replace the single "Codes / All codes"
#    description rows with an Operation table (Property/Operation/Value header,
#    then a code/regex/[0-9]{3} row), keeping the System URI rows as-is.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Include from ISO 3166-1 Codes")

# --- Sheet1 (Metadata): Version 0.1.2 -> 0.1.3 ---
$ws1.Range("B3").Value = "0.1.3"

# --- Sheet2: grow the header/content rows from 2 columns to 3 ---
# First, extend the formatting of the header row (A1) and data row (A2) across
# the new columns B and C so the new cells pick up the same styles.
$ws2.Range("A1").Copy()
$ws2.Range("B1:C1").PasteSpecial(-4122)   # xlPasteFormats

$ws2.Range("A2").Copy()
$ws2.Range("B2:C2").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# Now set the new cell values.
$ws2.Range("A1").Value = "Property"
$ws2.Range("B1").Value = "Operation"
$ws2.Range("C1").Value = "Value"

$ws2.Range("A2").Value = "code"
$ws2.Range("B2").Value = "regex"
$ws2.Range("C2").Value = "[0-9]{3}"

# Rows 3 and 4 (the blank spacer row and the System URI row) keep their
# existing content/style; nothing else to change there.
